# Add two new "summary" metric rows to the Metrics table:
#   - "# Residential Aged Care Molnupiravir Prescriptions" (before the existing
#     "...(Daily)" row)
#   - "# Residential Aged Care Paxlovid Prescriptions" (before the existing
#     "...(Daily)" row)
# and renumber the "Metric - Sort" column (D) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row right above the Molnupiravir (Daily) row (row 75),
# and another blank row above the Paxlovid (Daily) row, which after the
# first insert sits at row 78.
$ws.Rows("75:75").Insert()
$ws.Rows("78:78").Insert()

# New row 75: summary Molnupiravir metric
$ws.Range("A75").Value = "Treatments"
$ws.Range("B75").Value = 70
$ws.Range("C75").Value = "# Residential Aged Care Molnupiravir Prescriptions"
$ws.Range("D75").Value = 770
$ws.Range("F75").Value = "X"

# Existing Molnupiravir (Daily) row shifted down to row 76 - renumber sort col
$ws.Range("D76").Value = 780

# New row 78: summary Paxlovid metric
$ws.Range("A78").Value = "Treatments"
$ws.Range("B78").Value = 70
$ws.Range("C78").Value = "# Residential Aged Care Paxlovid Prescriptions"
$ws.Range("D78").Value = 800
$ws.Range("F78").Value = "X"

# Existing Paxlovid (Daily) row shifted down to row 79 - renumber sort col
$ws.Range("D79").Value = 810

# Resize the "Metrics" table to cover the two new rows (A1:F81 -> A1:F83)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F83"))

# Match the author's final selection/view state
$null = $ws.Range("F74:F78").Select()
